# Update the subreddit-corpus sheet:
#  - Row 2: new conv_id, clear question_text ([deleted] -> empty)
#  - Row 3: new conv_id, clear question_text ([deleted] -> empty)
#  - Row 4: new conv_id, clear question_text ([deleted] -> empty)
#  - Row 5: new conv_id (same as row4), res_id 1->2, question_text -> real text, clear answer_text (N/A -> empty)
#  - Row 6: new conv_id only (question/answer text stay "[deleted]"/"N/A")
#  - Row 7: new conv_id, clear question_text
#  - Row 8: new conv_id, clear question_text
#  - Rows 9-15: brand new rows appended to the corpus

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "30pk7l"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "N/A"

# --- Row 3 ---
$ws.Range("A3").Value = "58ccmp"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "N/A"

# --- Row 4 ---
$ws.Range("A4").Value = "39nogp"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "N/A"

# --- Row 5 ---
$ws.Range("A5").Value = "39nogp"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Sggvgcjcjvgc"
$ws.Range("D5").Value = ""

# --- Row 6 (conv_id only changes) ---
$ws.Range("A6").Value = "77swzz"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "[deleted]"
$ws.Range("D6").Value = "N/A"

# --- Row 7 ---
$ws.Range("A7").Value = "2s8qgt"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "N/A"

# --- Row 8 ---
$ws.Range("A8").Value = "63uwn9"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "N/A"

# --- Row 9 (new) ---
$ws.Range("A9").Value = "6uuo26"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "RAID THIS he hurt me today"
$ws.Range("D9").Value = "N/A"

$hockeyAwards = "After making all these I got a feel for what I liked and didn't, so juuuust for fun, Here are some awards - don't @ me.Bonus awards_____________Best Logo : San Jose SharksRunner up : Florida PanthersMost creative logo : Minnesota WildRunner up : Colorado AvalancheBest colors : Edmonton OilersRunner up : Arizona CoyotesWorst logo : Carolina HurricanesRunner up : Las Vegas Golden KnightsWorst colors : Las Vegas Golden KnightsRunner up : New Jersey DevilsBest alternate/retro logo : Anaheim DucksRunner up : Calgary FlamesMost creative alternate/retro logo : Arizona (Pheonix) CoyotesRunner up : Detroit Red WingsWorst alternate/retro logo : Ottawa SenatorsRunner up : Los Angeles Kings"

# --- Row 10 (new) ---
$ws.Range("A10").Value = "6x3ygl"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = $hockeyAwards
$ws.Range("D10").Value = "N/A"

# --- Row 11 (new) ---
$ws.Range("A11").Value = "6x3ygl"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "[deleted]"
$ws.Range("D11").Value = $hockeyAwards

# --- Row 12 (new) ---
$ws.Range("A12").Value = "6x3ygl"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "[deleted]"
$ws.Range("D12").Value = $hockeyAwards

# --- Row 13 (new) ---
$ws.Range("A13").Value = "6x3ygl"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "It’s-a-me, Mario, and I love lasagna!"
$ws.Range("D13").Value = "[deleted]"

# --- Row 14 (new) ---
$ws.Range("A14").Value = "3qcw3c"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = "N/A"

# --- Row 15 (new) ---
$ws.Range("A15").Value = "6eu2g1"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = "N/A"
